# Generate Report for Handoff
#
# - Status moves from "In Translation" to "Ready for handoff" everywhere it
#   appears (Overview!E2/F2, zh-cn!C2, de-de!C2 all share the same string).
# - The per-language "last generated" timestamps advance a few seconds:
#     Overview!G2            (Latest HO Xliff Generate Date)   03:00:29 -> 03:00:57
#     zh-cn!H2 (Latest Handoff Datetime)                       03:00:24 -> 03:00:52
#     de-de!H2 (Latest Handback DateTime, shares Overview's G2 string) -> 03:00:57
# - Columns E/F on Overview and column C on each language sheet (the
#   "zh-cn"/"de-de" / "Status" columns) widen a bit to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Timestamps --------------------------------------------------------------
$overview.Range("G2").Value = "2016-08-17 03:00:57"
$dede.Range("H2").Value     = "2016-08-17 03:00:57"
$zhcn.Range("H2").Value     = "2016-08-17 03:00:52"

# --- Column widths -----------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 16.33   # E
$overview.Columns.Item(6).ColumnWidth = 16.33   # F
$zhcn.Columns.Item(3).ColumnWidth     = 16.33   # C
$dede.Columns.Item(3).ColumnWidth     = 16.33   # C
